{"js": "// Update the worksheet date and every \"A\u00f7B=\" division-practice cell to the\n// new values from the commit. Each old value is unique in the document, so\n// an exact-text search & in-place replace is sufficient and keeps every\n// run's formatting (font/size) untouched.\nconst replacements = [\n  [\"2025-09-02 Tuesday\", \"2025-09-03 Wednesday\"],\n  [\"441\u00f79=\", \"519\u00f75=\"],\n  [\"270\u00f72=\", \"248\u00f76=\"],\n  [\"768\u00f75=\", \"441\u00f77=\"],\n  [\"367\u00f77=\", \"376\u00f75=\"],\n  [\"742\u00f76=\", \"381\u00f76=\"],\n  [\"270\u00f75=\", \"467\u00f79=\"],\n  [\"572\u00f74=\", \"609\u00f73=\"],\n  [\"239\u00f78=\", \"373\u00f79=\"],\n  [\"325\u00f78=\", \"170\u00f79=\"],\n  [\"342\u00f79=\", \"965\u00f79=\"],\n  [\"783\u00f73=\", \"691\u00f74=\"],\n  [\"746\u00f74=\", \"703\u00f75=\"],\n  [\"174\u00f76=\", \"413\u00f76=\"],\n  [\"711\u00f73=\", \"630\u00f76=\"],\n  [\"532\u00f72=\", \"875\u00f75=\"],\n  [\"367\u00f73=\", \"628\u00f75=\"],\n  [\"233\u00f74=\", \"601\u00f74=\"],\n  [\"660\u00f78=\", \"382\u00f78=\"],\n  [\"851\u00f77=\", \"745\u00f78=\"],\n  [\"242\u00f76=\", \"872\u00f74=\"],\n  [\"218\u00f78=\", \"348\u00f77=\"],\n  [\"932\u00f77=\", \"477\u00f78=\"],\n  [\"929\u00f77=\", \"436\u00f73=\"],\n  [\"129\u00f74=\", \"688\u00f72=\"],\n  [\"446\u00f74=\", \"896\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"A\u00f7B=\" division-practice cell to the\n# new values from the commit. Each old value is unique in the document, so\n# Find/Replace across the whole story is sufficient and preserves each\n# run's formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-09-02 Tuesday', '2025-09-03 Wednesday'),\n    @('441\u00f79=', '519\u00f75='),\n    @('270\u00f72=', '248\u00f76='),\n    @('768\u00f75=', '441\u00f77='),\n    @('367\u00f77=', '376\u00f75='),\n    @('742\u00f76=', '381\u00f76='),\n    @('270\u00f75=', '467\u00f79='),\n    @('572\u00f74=', '609\u00f73='),\n    @('239\u00f78=', '373\u00f79='),\n    @('325\u00f78=', '170\u00f79='),\n    @('342\u00f79=', '965\u00f79='),\n    @('783\u00f73=', '691\u00f74='),\n    @('746\u00f74=', '703\u00f75='),\n    @('174\u00f76=', '413\u00f76='),\n    @('711\u00f73=', '630\u00f76='),\n    @('532\u00f72=', '875\u00f75='),\n    @('367\u00f73=', '628\u00f75='),\n    @('233\u00f74=', '601\u00f74='),\n    @('660\u00f78=', '382\u00f78='),\n    @('851\u00f77=', '745\u00f78='),\n    @('242\u00f76=', '872\u00f74='),\n    @('218\u00f78=', '348\u00f77='),\n    @('932\u00f77=', '477\u00f78='),\n    @('929\u00f77=', '436\u00f73='),\n    @('129\u00f74=', '688\u00f72='),\n    @('446\u00f74=', '896\u00f74='),\n)\n\n# wdFindContinue = 1 (Wrap), wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n"}
